$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 4281
$ws.Range("F7").Value = 6296
$ws.Range("F11").Value = 29
$ws.Range("F12").Value = 9521
$ws.Range("F14").Value = 2572
$ws.Range("F16").Value = 2362
$ws.Range("F17").Value = 2573
$ws.Range("F19").Value = 263
$ws.Range("F20").Value = 2026
$ws.Range("F22").Value = 70
$ws.Range("F23").Value = 352
$ws.Range("F25").Value = 56
$ws.Range("F26").Value = 290
$ws.Range("F27").Value = 51
$ws.Range("F28").Value = 114
$ws.Range("F35").Value = 1614
$ws.Range("F36").Value = 2677
$ws.Range("F38").Value = 955
$ws.Range("F41").Value = 26

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 156
$ws.Range("F16").Value = 141

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 716
$ws.Range("F3").Value = 929

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 716
$ws.Range("F4").Value = 929
$ws.Range("F9").Value = 4281
$ws.Range("F10").Value = 4281
$ws.Range("F11").Value = 6296
$ws.Range("F14").Value = 29
$ws.Range("F15").Value = 9521
$ws.Range("F16").Value = 156
$ws.Range("F18").Value = 2572
$ws.Range("F20").Value = 2362
$ws.Range("F21").Value = 2573
$ws.Range("F23").Value = 263
$ws.Range("F24").Value = 2026
$ws.Range("F26").Value = 70
$ws.Range("F27").Value = 352
$ws.Range("F29").Value = 56
$ws.Range("F30").Value = 290
$ws.Range("F31").Value = 51
$ws.Range("F32").Value = 114
$ws.Range("F38").Value = 1614
$ws.Range("F40").Value = 2677
$ws.Range("F41").Value = 955
$ws.Range("F50").Value = 141
$ws.Range("F51").Value = 141
